$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EventData")
$ws.Activate()

# Row 3
$ws.Range("D3").Value = "eid event"
$ws.Range("C3").Value = "18"
$ws.Range("F3").Value = "No"

# Row 5
$ws.Range("C5").Value = "21"
$ws.Range("D5").Value = "sony proibar"
$ws.Range("E5").Value = "Success"
$ws.Range("F5").Value = "Yes"

$ws.Range("E9").Select()
